# Auto-generated script to apply row-shuffle + append 4 new rows
# matching commit: Atualizado por script em 05-11-2023 20:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowRange {
    param($rowNum, $values)
    $arr = New-Object 'object[,]' 1,17
    for ($i = 0; $i -lt $values.Count; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range("F" + $rowNum + ":V" + $rowNum).Value = $arr
}

$values_39 = @(
    'Villarreal B',
    1,
    'FC Cartagena SAD',
    2,
    2.37,
    '26/08/2023 20:42',
    2.11,
    '02/09/2023 20:56',
    3.33,
    '26/08/2023 20:42',
    3.47,
    '02/09/2023 20:52',
    3.15,
    '26/08/2023 20:42',
    3.75,
    '02/09/2023 20:59',
    'https://www.betexplorer.com/football/spain/laliga2/villarreal-fc-cartagena-sad/2FXXUSq4/'
)
Set-RowRange 39 $values_39

$values_40 = @(
    'Elche',
    1,
    'Racing Santander',
    1,
    1.65,
    '28/08/2023 23:42',
    1.92,
    '02/09/2023 20:57',
    3.82,
    '28/08/2023 23:42',
    3.51,
    '02/09/2023 20:57',
    5.66,
    '28/08/2023 23:42',
    4.46,
    '02/09/2023 20:57',
    'https://www.betexplorer.com/football/spain/laliga2/elche-racing-santander/CWkCDRET/'
)
Set-RowRange 40 $values_40

$values_51 = @(
    'Ferrol',
    2,
    'Villarreal B',
    2,
    2.22,
    '04/09/2023 11:12',
    1.88,
    '10/09/2023 16:08',
    3.13,
    '04/09/2023 11:12',
    3.53,
    '10/09/2023 16:08',
    3.71,
    '04/09/2023 11:12',
    4.63,
    '10/09/2023 16:08',
    'https://www.betexplorer.com/football/spain/laliga2/ferrol-villarreal/lzCGM4ip/'
)
Set-RowRange 51 $values_51

$values_52 = @(
    'Racing Santander',
    1,
    'Amorebieta',
    0,
    1.85,
    '04/09/2023 11:12',
    1.8,
    '10/09/2023 16:14',
    3.44,
    '04/09/2023 11:12',
    3.47,
    '10/09/2023 16:14',
    4.79,
    '04/09/2023 11:12',
    5.36,
    '10/09/2023 16:14',
    'https://www.betexplorer.com/football/spain/laliga2/racing-santander-amorebieta/AucBrnE3/'
)
Set-RowRange 52 $values_52

$values_53 = @(
    'FC Cartagena SAD',
    1,
    'Zaragoza',
    3,
    2.76,
    '03/09/2023 17:43',
    3.1,
    '10/09/2023 18:23',
    3.17,
    '03/09/2023 17:43',
    2.96,
    '10/09/2023 17:10',
    2.77,
    '03/09/2023 17:43',
    2.7,
    '10/09/2023 18:22',
    'https://www.betexplorer.com/football/spain/laliga2/fc-cartagena-sad-zaragoza/WfPlRpqT/'
)
Set-RowRange 53 $values_53

$values_54 = @(
    'Mirandes',
    4,
    'Andorra',
    3,
    2.85,
    '03/09/2023 20:12',
    2.88,
    '10/09/2023 18:23',
    3.16,
    '03/09/2023 20:12',
    3.12,
    '10/09/2023 17:08',
    2.75,
    '03/09/2023 20:12',
    2.79,
    '10/09/2023 18:27',
    'https://www.betexplorer.com/football/spain/laliga2/mirandes-fc-andorra/xUmOAStc/'
)
Set-RowRange 54 $values_54

$values_59 = @(
    'Valladolid',
    1,
    'FC Cartagena SAD',
    0,
    1.93,
    '13/09/2023 08:25',
    1.78,
    '16/09/2023 18:24',
    3.57,
    '13/09/2023 08:25',
    3.71,
    '16/09/2023 18:24',
    4.16,
    '13/09/2023 08:25',
    5.02,
    '16/09/2023 18:27',
    'https://www.betexplorer.com/football/spain/laliga2/valladolid-fc-cartagena-sad/lhXkfqDF/'
)
Set-RowRange 59 $values_59

$values_60 = @(
    'Alcorcon',
    0,
    'Levante',
    2,
    3.14,
    '11/09/2023 20:13',
    2.46,
    '16/09/2023 18:29',
    3.07,
    '11/09/2023 20:13',
    2.96,
    '16/09/2023 18:28',
    2.53,
    '11/09/2023 20:13',
    3.53,
    '16/09/2023 18:29',
    'https://www.betexplorer.com/football/spain/laliga2/alcorcon-levante/tUzfzOLq/'
)
Set-RowRange 60 $values_60

$values_72 = @(
    'Alcorcon',
    0,
    'Huesca',
    2,
    2.4,
    '18/09/2023 11:42',
    2.31,
    '23/09/2023 20:55',
    2.95,
    '18/09/2023 11:42',
    2.91,
    '23/09/2023 20:56',
    3.54,
    '18/09/2023 11:42',
    3.98,
    '23/09/2023 20:56',
    'https://www.betexplorer.com/football/spain/laliga2/alcorcon-huesca/2aDQmsjq/'
)
Set-RowRange 72 $values_72

$values_73 = @(
    'Racing Santander',
    2,
    'Albacete',
    1,
    2.56,
    '16/09/2023 20:12',
    2.66,
    '23/09/2023 20:57',
    3.12,
    '16/09/2023 20:12',
    3,
    '23/09/2023 20:57',
    3.12,
    '16/09/2023 20:12',
    3.15,
    '23/09/2023 20:57',
    'https://www.betexplorer.com/football/spain/laliga2/racing-santander-albacete/l6CO8rS8/'
)
Set-RowRange 73 $values_73

$values_77 = @(
    'Tenerife',
    1,
    'Espanyol',
    0,
    2.62,
    '24/09/2023 16:13',
    2.73,
    '25/09/2023 20:20',
    3.07,
    '24/09/2023 16:13',
    2.99,
    '25/09/2023 20:20',
    3.02,
    '24/09/2023 16:13',
    3.06,
    '25/09/2023 20:20',
    'https://www.betexplorer.com/football/spain/laliga2/tenerife-espanyol/xUT0kfRm/'
)
Set-RowRange 77 $values_77

$values_78 = @(
    'Ferrol',
    1,
    'Zaragoza',
    0,
    2.69,
    '18/09/2023 11:42',
    2.26,
    '25/09/2023 20:59',
    3.01,
    '18/09/2023 11:42',
    3.06,
    '25/09/2023 20:59',
    3,
    '18/09/2023 11:42',
    3.85,
    '25/09/2023 20:59',
    'https://www.betexplorer.com/football/spain/laliga2/ferrol-zaragoza/YcUyouz2/'
)
Set-RowRange 78 $values_78

$values_82 = @(
    'Albacete',
    3,
    'Andorra',
    1,
    1.96,
    '23/09/2023 20:12',
    1.93,
    '30/09/2023 20:41',
    3.52,
    '23/09/2023 20:12',
    3.44,
    '30/09/2023 20:59',
    4.07,
    '23/09/2023 20:12',
    4.55,
    '30/09/2023 20:54',
    'https://www.betexplorer.com/football/spain/laliga2/albacete-fc-andorra/Opjz8QQF/'
)
Set-RowRange 82 $values_82

$values_83 = @(
    'Villarreal B',
    2,
    'Alcorcon',
    2,
    2.51,
    '25/09/2023 11:42',
    2.59,
    '30/09/2023 20:40',
    3.15,
    '25/09/2023 11:42',
    3.22,
    '30/09/2023 20:46',
    3.18,
    '25/09/2023 11:42',
    3.02,
    '30/09/2023 20:41',
    'https://www.betexplorer.com/football/spain/laliga2/villarreal-alcorcon/44mr2XAg/'
)
Set-RowRange 83 $values_83

$values_85 = @(
    'Valladolid',
    3,
    'Burgos CF',
    0,
    1.81,
    '25/09/2023 11:42',
    1.87,
    '01/10/2023 18:26',
    3.61,
    '25/09/2023 11:42',
    3.43,
    '01/10/2023 18:26',
    4.86,
    '25/09/2023 11:42',
    4.95,
    '01/10/2023 18:26',
    'https://www.betexplorer.com/football/spain/laliga2/valladolid-burgos-cf/Gvin1iQa/'
)
Set-RowRange 85 $values_85

$values_86 = @(
    'Zaragoza',
    0,
    'Mirandes',
    1,
    1.88,
    '01/10/2023 16:13',
    1.88,
    '01/10/2023 16:13',
    3.38,
    '01/10/2023 16:13',
    3.38,
    '01/10/2023 16:13',
    4.97,
    '01/10/2023 16:13',
    4.97,
    '01/10/2023 16:13',
    'https://www.betexplorer.com/football/spain/laliga2/zaragoza-mirandes/b3ti0Bu6/'
)
Set-RowRange 86 $values_86

$values_87 = @(
    'Amorebieta',
    0,
    'FC Cartagena SAD',
    0,
    2.38,
    '25/09/2023 11:42',
    2.72,
    '01/10/2023 18:29',
    3.24,
    '25/09/2023 11:42',
    3.09,
    '01/10/2023 18:06',
    3.23,
    '25/09/2023 11:42',
    2.98,
    '01/10/2023 18:29',
    'https://www.betexplorer.com/football/spain/laliga2/amorebieta-fc-cartagena-sad/2yRCnhB6/'
)
Set-RowRange 87 $values_87

$values_90 = @(
    'Alcorcon',
    1,
    'Albacete',
    2,
    2.82,
    '01/10/2023 04:42',
    3.14,
    '03/10/2023 18:59',
    3.11,
    '01/10/2023 04:42',
    3.04,
    '03/10/2023 18:42',
    2.77,
    '01/10/2023 04:42',
    2.63,
    '03/10/2023 18:59',
    'https://www.betexplorer.com/football/spain/laliga2/alcorcon-albacete/OpEDHCAt/'
)
Set-RowRange 90 $values_90

$values_91 = @(
    'Levante',
    1,
    'Villarreal B',
    1,
    1.56,
    '30/09/2023 20:12',
    1.66,
    '03/10/2023 18:54',
    4.32,
    '30/09/2023 20:12',
    4.17,
    '03/10/2023 18:54',
    6.06,
    '30/09/2023 20:12',
    5.31,
    '03/10/2023 18:54',
    'https://www.betexplorer.com/football/spain/laliga2/levante-villarreal/nPMoMh2P/'
)
Set-RowRange 91 $values_91

$values_93 = @(
    'R. Oviedo',
    1,
    'Huesca',
    0,
    2.17,
    '01/10/2023 20:12',
    1.85,
    '04/10/2023 18:54',
    3.1,
    '01/10/2023 20:12',
    3.15,
    '04/10/2023 18:55',
    3.92,
    '01/10/2023 20:12',
    5.82,
    '04/10/2023 18:55',
    'https://www.betexplorer.com/football/spain/laliga2/r-oviedo-huesca/2BsMDU95/'
)
Set-RowRange 93 $values_93

$values_94 = @(
    'Eldense',
    0,
    'Valladolid',
    1,
    2.92,
    '01/10/2023 17:43',
    3.61,
    '04/10/2023 18:59',
    3.21,
    '01/10/2023 17:43',
    3.32,
    '04/10/2023 18:17',
    2.6,
    '01/10/2023 17:43',
    2.23,
    '04/10/2023 18:55',
    'https://www.betexplorer.com/football/spain/laliga2/eldense-valladolid/pKr3c9PO/'
)
Set-RowRange 94 $values_94

$values_95 = @(
    'Gijon',
    2,
    'Elche',
    0,
    2.47,
    '01/10/2023 20:12',
    2.39,
    '04/10/2023 21:26',
    3.29,
    '01/10/2023 20:12',
    3.09,
    '04/10/2023 21:26',
    3.1,
    '01/10/2023 20:12',
    3.5,
    '04/10/2023 21:26',
    'https://www.betexplorer.com/football/spain/laliga2/gijon-elche/OKrIEAfa/'
)
Set-RowRange 95 $values_95

$values_96 = @(
    'Burgos CF',
    1,
    'Leganes',
    0,
    2.23,
    '01/10/2023 17:43',
    2.29,
    '04/10/2023 21:28',
    3.14,
    '01/10/2023 17:43',
    2.9,
    '04/10/2023 21:28',
    3.68,
    '01/10/2023 17:43',
    4.04,
    '04/10/2023 21:29',
    'https://www.betexplorer.com/football/spain/laliga2/burgos-cf-leganes/xfueaVfC/'
)
Set-RowRange 96 $values_96

$values_99 = @(
    'FC Cartagena SAD',
    0,
    'Espanyol',
    2,
    4.05,
    '02/10/2023 20:12',
    3.73,
    '05/10/2023 21:28',
    3.6,
    '02/10/2023 20:12',
    3.45,
    '05/10/2023 21:28',
    1.94,
    '02/10/2023 20:12',
    2.13,
    '05/10/2023 21:26',
    'https://www.betexplorer.com/football/spain/laliga2/fc-cartagena-sad-espanyol/6TqabkAI/'
)
Set-RowRange 99 $values_99

$values_100 = @(
    'Ferrol',
    1,
    'Amorebieta',
    0,
    1.74,
    '02/10/2023 20:12',
    1.71,
    '05/10/2023 21:21',
    3.61,
    '02/10/2023 20:12',
    3.58,
    '05/10/2023 21:28',
    5.24,
    '02/10/2023 20:12',
    5.95,
    '05/10/2023 21:28',
    'https://www.betexplorer.com/football/spain/laliga2/ferrol-amorebieta/I7DHGWPn/'
)
Set-RowRange 100 $values_100

$values_114 = @(
    'Gijon',
    2,
    'Zaragoza',
    2,
    2.32,
    '08/10/2023 20:12',
    1.99,
    '14/10/2023 18:21',
    3.03,
    '08/10/2023 20:12',
    3.31,
    '14/10/2023 18:29',
    3.59,
    '08/10/2023 20:12',
    4.49,
    '14/10/2023 18:29',
    'https://www.betexplorer.com/football/spain/laliga2/gijon-zaragoza/4bSYumgb/'
)
Set-RowRange 114 $values_114

$values_115 = @(
    'FC Cartagena SAD',
    2,
    'Racing Santander',
    3,
    2.31,
    '08/10/2023 17:43',
    2.39,
    '14/10/2023 18:22',
    3.31,
    '08/10/2023 17:43',
    3.1,
    '14/10/2023 18:22',
    3.38,
    '08/10/2023 17:43',
    3.48,
    '14/10/2023 18:26',
    'https://www.betexplorer.com/football/spain/laliga2/fc-cartagena-sad-racing-santander/vHbGN5wU/'
)
Set-RowRange 115 $values_115

$values_120 = @(
    'Eibar',
    1,
    'Huesca',
    1,
    1.57,
    '08/10/2023 17:43',
    1.59,
    '15/10/2023 18:22',
    3.87,
    '08/10/2023 17:43',
    3.7,
    '15/10/2023 18:22',
    6.75,
    '08/10/2023 17:43',
    7.53,
    '15/10/2023 18:24',
    'https://www.betexplorer.com/football/spain/laliga2/eibar-huesca/GxjrJPoo/'
)
Set-RowRange 120 $values_120

$values_121 = @(
    'Leganes',
    6,
    'Amorebieta',
    0,
    1.7,
    '09/10/2023 11:11',
    1.74,
    '15/10/2023 18:23',
    3.57,
    '09/10/2023 11:11',
    3.63,
    '15/10/2023 18:23',
    5.98,
    '09/10/2023 11:11',
    5.58,
    '15/10/2023 18:23',
    'https://www.betexplorer.com/football/spain/laliga2/leganes-amorebieta/QDGsvRNA/'
)
Set-RowRange 121 $values_121

$values_147 = @(
    'Valladolid',
    2,
    'Tenerife',
    0,
    1.99,
    '28/10/2023 21:12',
    2.26,
    '04/11/2023 18:24',
    3.32,
    '28/10/2023 21:12',
    3.03,
    '04/11/2023 18:24',
    4.39,
    '28/10/2023 21:12',
    3.91,
    '04/11/2023 18:24',
    'https://www.betexplorer.com/football/spain/laliga2/valladolid-tenerife/Gb3CsLo0/'
)
Set-RowRange 147 $values_147

$values_148 = @(
    'Levante',
    2,
    'Mirandes',
    2,
    1.65,
    '29/10/2023 14:12',
    1.9,
    '04/11/2023 18:23',
    4.24,
    '29/10/2023 14:12',
    3.53,
    '04/11/2023 18:23',
    4.95,
    '29/10/2023 14:12',
    4.52,
    '04/11/2023 18:23',
    'https://www.betexplorer.com/football/spain/laliga2/levante-mirandes/6PTSkqvK/'
)
Set-RowRange 148 $values_148

# Append new rows 151-154
$ws.Range("A150").Copy()
$ws.Range("A151:A154").PasteSpecial(-4122)
$ws.Range("E150").Copy()
$ws.Range("E151:E154").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A151").Value = 150
$ws.Range("B151").Value = 'spain'
$ws.Range("C151").Value = 'laliga2'
$ws.Range("D151").Value = '2023-2024'
$ws.Range("E151").Value = 45235.67708333334
$values_151 = @(
    'Amorebieta',
    0,
    'Huesca',
    1,
    2.65,
    '29/10/2023 16:43',
    2.92,
    '05/11/2023 16:13',
    3.01,
    '29/10/2023 16:43',
    2.94,
    '05/11/2023 16:07',
    3.1,
    '29/10/2023 16:43',
    2.9,
    '05/11/2023 16:13',
    'https://www.betexplorer.com/football/spain/laliga2/amorebieta-huesca/EVtsSunD/'
)
Set-RowRange 151 $values_151

$ws.Range("A152").Value = 151
$ws.Range("B152").Value = 'spain'
$ws.Range("C152").Value = 'laliga2'
$ws.Range("D152").Value = '2023-2024'
$ws.Range("E152").Value = 45235.77083333334
$values_152 = @(
    'FC Cartagena SAD',
    0,
    'Leganes',
    3,
    2.38,
    '29/10/2023 18:42',
    2.98,
    '05/11/2023 18:22',
    3.18,
    '29/10/2023 18:42',
    3.03,
    '05/11/2023 18:24',
    3.29,
    '29/10/2023 18:42',
    2.76,
    '05/11/2023 18:29',
    'https://www.betexplorer.com/football/spain/laliga2/fc-cartagena-sad-leganes/SYUOjPOD/'
)
Set-RowRange 152 $values_152

$ws.Range("A153").Value = 152
$ws.Range("B153").Value = 'spain'
$ws.Range("C153").Value = 'laliga2'
$ws.Range("D153").Value = '2023-2024'
$ws.Range("E153").Value = 45235.77083333334
$values_153 = @(
    'Alcorcon',
    3,
    'Racing Santander',
    1,
    2.21,
    '29/10/2023 21:12',
    2.63,
    '05/11/2023 18:29',
    3.25,
    '29/10/2023 21:12',
    3.22,
    '05/11/2023 18:28',
    3.58,
    '29/10/2023 21:12',
    2.97,
    '05/11/2023 18:29',
    'https://www.betexplorer.com/football/spain/laliga2/alcorcon-racing-santander/2ewwTLW6/'
)
Set-RowRange 153 $values_153

$ws.Range("A154").Value = 153
$ws.Range("B154").Value = 'spain'
$ws.Range("C154").Value = 'laliga2'
$ws.Range("D154").Value = '2023-2024'
$ws.Range("E154").Value = 45235.875
$values_154 = @(
    'Eldense',
    2,
    'Burgos CF',
    0,
    2.15,
    '30/10/2023 21:13',
    2.27,
    '05/11/2023 20:36',
    3.19,
    '30/10/2023 21:13',
    3.04,
    '05/11/2023 20:50',
    3.83,
    '30/10/2023 21:13',
    3.87,
    '05/11/2023 20:50',
    'https://www.betexplorer.com/football/spain/laliga2/eldense-burgos-cf/O0OXl3gQ/'
)
Set-RowRange 154 $values_154
